$d = $word.ActiveDocument

# Locate the "Ver no Jupiter ..." paragraph (start of the block to remove).
$rStart = $d.Content
$foundStart = $rStart.Find.Execute(
    "Ver no Jupiter Salvar em pdf Salvar em docx",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundStart) {
    throw "Could not find the 'Ver no Jupiter ...' paragraph"
}

# Locate the "© 2020 ..." paragraph (end of the text to remove).
$rEnd = $d.Content
$foundEnd = $rEnd.Find.Execute(
    "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundEnd) {
    throw "Could not find the '© 2020 ...' paragraph"
}

# Expand each found range to its whole paragraph (so the paragraph marks are
# included in the deletion, not just the run text).
$startPara = $rStart.Paragraphs(1)
$endPara = $rEnd.Paragraphs(1)

# The empty paragraph right after the "© 2020 ..." paragraph is also removed
# by the edit (it collapses into the blank paragraph that already follows
# "LOB1018: ..."), so grab it too.
$afterEnd = $d.Range($endPara.Range.End, $endPara.Range.End)
$trailingBlankPara = $afterEnd.Paragraphs(1)

$deleteRange = $d.Range($startPara.Range.Start, $trailingBlankPara.Range.End)
$deleteRange.Delete()
